# Apply the "Hjemme passive tweaks lichtwark deleted values" edit to the
# single data sheet (Ark1): row 1 gets new group-size header values, and a
# couple of cells in rows 2-3 get recalculated (some cells are cleared
# entirely, others get new numeric values, and D3 gains a new value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 (header / N values)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON)
$ws.Range("B2").Value = -8.9276162598221021
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 6.1040602419814824
$ws.Range("E2").Value = 6.5460661122605472

# Row 3 (STR)
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 16.67413550689594
$ws.Range("D3").Value = 22.123273101918144
$ws.Range("E3").Value = 2.3701403018050841

# Reflect the narrower selection left in the saved file
$ws.Range("B1:E3").Select()
